$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log Iter 1-10")
$ws.Range("Z100").Value = "test"
"done"
